$p = $ppt.ActivePresentation

# --- 1) Table style swap -------------------------------------------------
# Three slides (14, 15, 16) each contain a single table whose
# <a:tableStyleId> changes from {1592EB1B-8B87-4154-8588-1D75798A1A00} to
# {032BFA48-6943-453E-9126-F424D235407C}.
$newTableStyle = "{032BFA48-6943-453E-9126-F424D235407C}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Theme colour-scheme swap ----------------------------------------
# The deck's two theme parts had their contents swapped: the slide
# master's theme ("Integral" / Red Violet) and the notes master's theme
# ("Office Theme") traded places. The accessible surface for editing
# theme colours is Slide.ThemeColorScheme, which always resolves to the
# slide master's theme palette, so recolour it to the target (Office
# Theme) palette.
function Set-RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
$tcs.Colors(1).RGB  = Set-RGB 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = Set-RGB 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = Set-RGB 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = Set-RGB 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = Set-RGB 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = Set-RGB 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = Set-RGB 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = Set-RGB 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = Set-RGB 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = Set-RGB 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = Set-RGB 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = Set-RGB 0x95 0x4F 0x72   # folHlink
